$d = $word.ActiveDocument

# --- 1) FACULTY line: split the long underscore run and insert the
#        underlined faculty name in the middle of it ---
$rng = $d.Content
$rng.Find.Execute("ФАКУЛЬТЕТ ___________________________________________________________________", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = "ФАКУЛЬТЕТ ________"

$facName = $rng.Duplicate
$facName.Collapse(0)
$facName.InsertAfter("Информатика и системы управления")
$facName.Font.Underline = 1

$facTail = $d.Range($facName.End, $facName.End)
$facTail.InsertAfter("____________________________")
$facTail.Font.Size = 12

# --- 2) KAFEDRA line: split the long underscore run (the one carrying
#        the iCs property) and insert the underlined department name ---
$rng2 = $d.Content
$rng2.Find.Execute("_____________________________________________________________________", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.Text = "_________"

$depName = $rng2.Duplicate
$depName.Collapse(0)
$depName.InsertAfter("Программное обеспечение ЭВМ и информационные техлоногии")
$depName.Font.Italic = $depName.Font.Italic
$depName.Font.Underline = 1

$depTail = $d.Range($depName.End, $depName.End)
$depTail.InsertAfter("_____")
$depTail.Font.Size = 12

# --- 3) Drop the stray explicit en-US language override on the lone
#        "5" run inside "Группа__ИУ75__6Б__________" without merging it
#        into its neighbouring runs ---
$rng3 = $d.Content
$rng3.Find.Execute("Группа__ИУ7-5", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$digitRng = $d.Range($rng3.End - 1, $rng3.End)
$digitRng.Delete()
$digitIns = $d.Range($rng3.End - 1, $rng3.End - 1)
$digitIns.InsertAfter("5")
$digitFix = $d.Range($rng3.End - 1, $rng3.End)
$digitFix.Font.Bold = 1
$digitFix.Font.Bold = 0

Write-Output "done"
